$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the device count used in the run-time estimates (123 -> 117 apps),
# and bump the "apps/dev" split for row 9 (12 -> 13) to match.
$ws.Range("B7").Value = 117
$ws.Range("B9").Value = 117
$ws.Range("C9").Value = 13
$ws.Range("B11").Value = 117
$ws.Range("B14").Value = 117
$ws.Range("B16").Value = 117
$ws.Range("B18").Value = 117
$ws.Range("B20").Value = 117
$ws.Range("B22").Value = 117
$ws.Range("B24").Value = 117
$ws.Range("B26").Value = 117
$ws.Range("B28").Value = 117

# Move the active selection to B17, matching the author's saved cursor position.
$ws.Range("B17").Select()
